# top_videos_scored: refresh the scraped stats, add a human-readable
# "duration_str" column (inserted before video_id), re-rank, and drop one
# stale row (the two bottom rows collapse into one updated row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before H ("video_id") to hold duration_str; this
# shifts video_id..rank one column to the right (H..R -> I..S).
$ws.Columns("H").Insert()

# Drop the last row outright - its data now lives in (the refreshed) row 5.
$ws.Rows("6").Delete()

# Row 1
$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "channel"
$ws.Range("C1").Value = "published"
$ws.Range("D1").Value = "views"
$ws.Range("E1").Value = "likes"
$ws.Range("F1").Value = "comments"
$ws.Range("G1").Value = "duration_minutes"
$ws.Range("H1").Value = "duration_str"
$ws.Range("I1").Value = "video_id"
$ws.Range("J1").Value = "url"
$ws.Range("K1").Value = "likes_per_view"
$ws.Range("L1").Value = "comments_per_minute"
$ws.Range("M1").Value = "views_per_day"
$ws.Range("N1").Value = "norm_likes_per_view"
$ws.Range("O1").Value = "norm_comments_per_minute"
$ws.Range("P1").Value = "norm_views_per_day"
$ws.Range("Q1").Value = "norm_views"
$ws.Range("R1").Value = "final_score"
$ws.Range("S1").Value = "rank"

# Row 2
$ws.Range("A2").Value = "Algorithmic Trading Using Python - Full Course"
$ws.Range("B2").Value = "freeCodeCamp.org"
$ws.Range("C2").Value = 44169.72681712963
$ws.Range("D2").Value = 2940151
$ws.Range("E2").Value = 75970
$ws.Range("F2").Value = 2245
$ws.Range("G2").Value = 273.05
$ws.Range("H2").Value = "4:33:03"
$ws.Range("I2").Value = "xfzGZB4HhEE"
$ws.Range("J2").Value = "https://www.youtube.com/watch?v=xfzGZB4HhEE"
$ws.Range("K2").Value = 0.02583880895913169
$ws.Range("L2").Value = 8.218927329306242
$ws.Range("M2").Value = 1814.908024691358
$ws.Range("N2").Value = 0.350782384575447
$ws.Range("O2").Value = 0.03495405875667099
$ws.Range("P2").Value = 0.7974064119635269
$ws.Range("Q2").Value = 0.9999999999999997
$ws.Range("R2").Value = 5.514474507130263
$ws.Range("S2").Value = 1

# Row 3
$ws.Range("A3").Value = "Algorithmic Trading – Machine Learning & Quant Strategies Course with Python"
$ws.Range("B3").Value = "freeCodeCamp.org"
$ws.Range("C3").Value = 45225.52559027778
$ws.Range("D3").Value = 1003867
$ws.Range("E3").Value = 16878
$ws.Range("F3").Value = 374
$ws.Range("G3").Value = 179.33
$ws.Range("H3").Value = "2:59:20"
$ws.Range("I3").Value = "9Y3yaoi9rUQ"
$ws.Range("J3").Value = "https://www.youtube.com/watch?v=9Y3yaoi9rUQ"
$ws.Range("K3").Value = 0.01681298419013674
$ws.Range("L3").Value = 2.084378309089896
$ws.Range("M3").Value = 1779.906028368794
$ws.Range("N3").Value = 0.2282496339275403
$ws.Range("O3").Value = 0.008864597406436569
$ws.Range("P3").Value = 0.7816251143145888
$ws.Range("Q3").Value = 0.3414006720401305
$ws.Range("R3").Value = 3.730154783619521
$ws.Range("S3").Value = 2

# Row 4
$ws.Range("A4").Value = "He Makes a Living Algo Trading in Forex - Scott Welsh | Trader Interview"
$ws.Range("B4").Value = "Etienne Crete - Desire To TRADE"
$ws.Range("C4").Value = 44731.45846064815
$ws.Range("D4").Value = 48993
$ws.Range("E4").Value = 1100
$ws.Range("F4").Value = 68
$ws.Range("G4").Value = 34.47
$ws.Range("H4").Value = "34:28"
$ws.Range("I4").Value = "TKVE6DL7ubU"
$ws.Range("J4").Value = "https://www.youtube.com/watch?v=TKVE6DL7ubU"
$ws.Range("K4").Value = 0.02245218704712918
$ws.Range("L4").Value = 1.967023430720278
$ws.Range("M4").Value = 46.30718336483932
$ws.Range("N4").Value = 0.3048062983004695
$ws.Range("O4").Value = 0.008365501946705834
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0.0166139286252429
$ws.Range("R4").Value = 0.9643777560453058
$ws.Range("S4").Value = 3

# Row 5
$ws.Range("A5").Value = "7 Algo Trading Strategies (Backtest And Rules)"
$ws.Range("B5").Value = "Quantified Strategies"
$ws.Range("C5").Value = 45194.54184027778
$ws.Range("D5").Value = 31246
$ws.Range("E5").Value = 505
$ws.Range("F5").Value = 12
$ws.Range("G5").Value = 9.42
$ws.Range("H5").Value = "9:25"
$ws.Range("I5").Value = "NojfYk31_xI"
$ws.Range("J5").Value = "https://www.youtube.com/watch?v=NojfYk31_xI"
$ws.Range("K5").Value = 0.01616206874479933
$ws.Range("L5").Value = 1.260504201680672
$ws.Range("M5").Value = 52.51428571428571
$ws.Range("N5").Value = 0.2194129389996283
$ws.Range("O5").Value = 0.005360764995630636
$ws.Range("P5").Value = 0.002798586938051013
$ws.Range("Q5").Value = 0.01057754022699976
$ws.Range("R5").Value = 0.6985111882582986
$ws.Range("S5").Value = 4
